$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 64

$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "01/27/2026"
$ws.Cells.Item($row, 2).Value = 11394.85
$ws.Cells.Item($row, 3).Value = 0.2427946360339521
$ws.Cells.Item($row, 4).Value = 0.7572053639660479
$ws.Cells.Item($row, 5).Value = -215.67
$ws.Cells.Item($row, 6).Value = -29.34
$ws.Cells.Item($row, 7).Value = -22163.52
$ws.Cells.Item($row, 8).Value = -71.98
$ws.Cells.Item($row, 9).Value = -486.24
$ws.Cells.Item($row, 10).Value = -14.95
$ws.Cells.Item($row, 11).Value = -22649.76
$ws.Cells.Item($row, 12).Value = -66.53
